# Update description strings in the "Items" and "Items - Formatted" sheets
# to match the new test-case data (per commit "Update on Export.xaml Test Case").

$wb = $excel.ActiveWorkbook

$updates = @{
    "A2"  = "green onion Pancakes MMt4AF (1)"
    "A3"  = "Pan Fried Leek Dumplings #8T (2)"
    "A4"  = "Pork Xiao Long Bao(10) Ñ#P]]]NAÈ#Q(10)"
    "A5"  = "Q-BA0) (5) MEÀE (5)"
    "A6"  = "Chicken potstickers RÈP]]$/##5(6)"
    "A7"  = "Tomato Mushroom Steamed dumpli ₪ (6)"
    "A8"  = "Zucchini shrimp dumplings j JJJ#tl6/5"
    "A9"  = "beef stew nodle soup (Non Spicy 25+ØJ(T#)"
    "A10" = "dandan noodle #/m"
    "A11" = "banana naan bread TATRAI"
    "A12" = "house made plum juice"
}

foreach ($sheetName in @("Items", "Items - Formatted")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
